$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Split F2's text into two rich-text runs -----------------------------
# Original value: "50659,50660,50661"  ->  new value: "50659,50660,50661,50658"
# First run keeps the cell's existing font (sz 8, Arial); second (new) run
# is appended with its own (new) font entry.
$ws.Range("F2").Value2 = "50659,50660,50661,50658"

# Second run = the appended "50658" (characters 19-23, 1-based).
$ch2 = $ws.Range("F2").Characters(19, 5)
$ch2.Font.ColorIndex = -4105

# --- 2. Selection moves from B8 to F3 ---------------------------------------
$ws.Range("F3").Select()

# --- 3. Column layout: the single E:F (width 12.1) column group is replaced
#        by two separate, wider columns -------------------------------------
$ws.Columns("E").ColumnWidth = 19.5
$ws.Columns("F").ColumnWidth = 25.3

# --- 4. Register a dedicated (sz 8 / Arial) font in the style table --------
# Applying the font directly to a scratch cell far outside the used range
# (then deleting that row so no trace of it remains) forces the workbook to
# keep the font definition without leaving stray cells/rows behind.
$ws.Range("A100").Value2 = "temp"
$ws.Range("A100").Font.Size = 8
$ws.Range("A100").Font.Name = "Arial"
$ws.Range("A100").EntireRow.Delete()
